$wb = $excel.ActiveWorkbook

# --- CompanyUpdates sheet: flip the Subject/Client potential-round-trip flags ---
# (Source - Engagement stays the same; the Yes/No values swap columns)
$wsCompanyUpdates = $wb.Worksheets.Item("CompanyUpdates")
$wsCompanyUpdates.Range("A2:A5").Value = "No"
$wsCompanyUpdates.Range("C2:C5").Value = "Yes"
$wsCompanyUpdates.Range("C11").Select()

# --- AddOpportunity sheet: replace the round-trip example company pair and ---
# --- clear out the stale example rows, then make this the active tab       ---
$wsAddOpportunity = $wb.Worksheets.Item("AddOpportunity")
$wsAddOpportunity.Activate()
$wsAddOpportunity.Range("B2").Value = "Grupo Volum"
$wsAddOpportunity.Range("A2").Value = "GPF Capital"
$wsAddOpportunity.Range("A3:B5").ClearContents()
$wsAddOpportunity.Range("C17").Select()
